$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New period row (17-10-28 ~ 17-10-29) appended after the existing A1:N35
# data block. Every column in this report stores its value as literal text
# (shared string), even columns that look numeric/percentage (e.g. "81",
# "100.00%"), so plain non-numeric-looking values are written directly while
# digit-/percent-looking ones are forced to Text first -- otherwise Excel
# would silently convert them into a real number/percentage -- and the
# cell's style is then put back to "Normal" so it keeps the same (absent)
# per-cell style as its neighbours in the existing rows.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("A36").Value = '17-10-28 ~ 17-10-29'

Set-TextValue $ws.Range("B36") '81'
Set-TextValue $ws.Range("C36") '81'
Set-TextValue $ws.Range("D36") '100.00%'
Set-TextValue $ws.Range("E36") '0'
Set-TextValue $ws.Range("F36") '81'
Set-TextValue $ws.Range("G36") '100.00%'
Set-TextValue $ws.Range("H36") '0'
Set-TextValue $ws.Range("I36") '81'
Set-TextValue $ws.Range("J36") '0'
Set-TextValue $ws.Range("K36") '0.00%'

$ws.Range("L36").Value = '[]'
$ws.Range("L36").Style = "Normal"
$ws.Range("M36").Value = '[]'
$ws.Range("N36").Value = '[''7ff9010202000022'', ''7ff9010202000024'', ''7ff9010202000025'', ''7ff9010202000026'', ''7ff9010202000027'', ''7ff9010202000028'', ''7ff9010202000030'', ''7ff9010202000031'', ''7ff9010202000032'', ''7ff9010202000034'', ''7ff9010202000035'', ''7ff9010202000036'', ''7ff9010202000037'', ''7ff9010202000038'', ''7ff9010202000039'', ''7ff9010202000040'', ''7ff9010202000042'', ''7ff9010202000043'', ''7ff9010202000045'', ''7ff9010202000049'', ''7ff9010202000050'', ''7ff9010202000053'', ''7ff9010202000055'', ''7ff9010202000056'', ''7ff9010202000057'', ''7ff9010202000058'', ''7ff9010202000059'', ''7ff9010202000060'', ''7ff9010202000061'', ''7ff9010202000064'', ''7ff9010202000067'', ''7ff9010202000069'', ''7ff9010202000073'', ''7ff9010202000074'', ''7ff9010202000081'', ''7ff9010202000082'', ''7ff9010202000083'', ''7ff9010202000084'', ''7ff9010202000085'', ''7ff9010202000087'', ''7ff9010202000088'', ''7ff9010202000090'', ''7ff9010202000092'', ''7ff9010202000095'', ''7ff9010202000097'', ''7ff9010202000098'', ''7ff9010202000102'', ''7ff9010202000104'', ''7ff9010202000106'', ''7ff9010202000107'', ''7ff9010202000108'', ''7ff9010202000112'', ''7ff9010202000113'', ''7ff9010202000116'', ''7ff9010202000117'', ''7ff9010202000118'', ''7ff9010202000121'', ''7ff9010202000123'', ''7ff9010202000126'', ''7ff9010202000127'', ''7ff9010202000129'', ''7ff9010202000131'', ''7ff9010202000132'', ''7ff9010202000134'', ''7ff9010202000136'', ''7ff9010202000137'', ''7ff9010202000138'', ''7ff9010202000139'', ''7ff9010202000141'', ''7ff9010202000146'', ''7ff9010202000147'', ''7ff9010202000151'', ''7ff9010202000153'', ''7ff9010202000156'', ''7ff9010202000161'', ''7ff9010202000162'', ''7ff9010202000163'', ''7ff9010202000164'', ''7ff9010202000165'', ''7ff9010202000166'', ''7ff9010202000167'']'
